$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H32").Value = 923.3077
$ws.Range("I32").Value = 1100
$ws.Range("J32").Value = 870.3
$ws.Range("K32").Value = 1100
$ws.Range("L32").Value = 870.3
$ws.Range("M32").Value = -774
$ws.Range("N32").Value = -1522.3
$ws.Range("H51").Value = 3082.8333
$ws.Range("I51").Value = 2512.375
$ws.Range("J51").Value = 3539.2
$ws.Range("K51").Value = 2512.375
$ws.Range("L51").Value = 3539.2
$ws.Range("M51").Value = -2028.375
$ws.Range("N51").Value = -4507.2
$ws.Range("H74").Value = 3341.625
$ws.Range("I74").Value = 3183.25
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 3183.25
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -2247.25
$ws.Range("N74").Value = -5372
$ws.Range("H77").Value = 3341.625
$ws.Range("I77").Value = 3183.25
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 15916.25
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -11236.25
$ws.Range("N77").Value = -26860
$ws.Range("H139").Value = 30173.334
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 30173.334
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 30173.334
$ws.Range("N139").Value = -40453.334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 200
$ws.Range("I10").Value = 200
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 200
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -30
$ws.Range("N10").ClearContents()
$ws.Range("H32").Value = 9540.895500000001
$ws.Range("I32").Value = 10943.865
$ws.Range("J32").Value = 4821.8184
$ws.Range("K32").Value = 10943.865
$ws.Range("L32").Value = 4821.8184
$ws.Range("M32").Value = -10656.865
$ws.Range("N32").Value = -5395.8184
$ws.Range("H45").Value = 2200
$ws.Range("I45").Value = 2178
$ws.Range("J45").Value = 2266
$ws.Range("K45").Value = 2178
$ws.Range("L45").Value = 2266
$ws.Range("M45").Value = -1801
$ws.Range("N45").Value = -3020
$ws.Range("H61").Value = 14709164
$ws.Range("I61").Value = 16669686
$ws.Range("J61").Value = 5250
$ws.Range("K61").Value = 16669686
$ws.Range("L61").Value = 5250
$ws.Range("M61").Value = -16669474
$ws.Range("N61").Value = -5674
$ws.Range("H74").Value = 17860334
$ws.Range("I74").Value = 31252112
$ws.Range("J74").Value = 4632.1665
$ws.Range("K74").Value = 31252112
$ws.Range("L74").Value = 4632.1665
$ws.Range("M74").Value = -31251238
$ws.Range("N74").Value = -6380.1665
$ws.Range("H77").Value = 17860334
$ws.Range("I77").Value = 31252112
$ws.Range("J77").Value = 4632.1665
$ws.Range("K77").Value = 156260560
$ws.Range("L77").Value = 23160.8325
$ws.Range("M77").Value = -156256192
$ws.Range("N77").Value = -31896.8325
$ws.Range("H102").Value = 870
$ws.Range("I102").Value = 1110
$ws.Range("J102").Value = 750
$ws.Range("K102").Value = 1110
$ws.Range("L102").Value = 750
$ws.Range("M102").Value = 512
$ws.Range("N102").Value = -3994
$ws.Range("H110").Value = 1237.2
$ws.Range("I110").Value = 1128.6364
$ws.Range("J110").Value = 2033.3334
$ws.Range("K110").Value = 1128.6364
$ws.Range("L110").Value = 2033.3334
$ws.Range("M110").Value = 916.3635999999999
$ws.Range("N110").Value = -6123.3334
$ws.Range("H136").Value = 14709164
$ws.Range("I136").Value = 16669686
$ws.Range("J136").Value = 5250
$ws.Range("K136").Value = 50009058
$ws.Range("L136").Value = 15750
$ws.Range("M136").Value = -50006508
$ws.Range("N136").Value = -20850

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4032.2
$ws.Range("I105").Value = 2685.8667
$ws.Range("J105").Value = 4840
$ws.Range("K105").Value = 2685.8667
$ws.Range("L105").Value = 4840
$ws.Range("M105").Value = -938.8667
$ws.Range("N105").Value = -8334

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 47328.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 47328.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 47328.5
$ws.Range("N43").Value = -47696.5
$ws.Range("H58").Value = 3284.2666
$ws.Range("I58").Value = 1244.875
$ws.Range("J58").Value = 5615
$ws.Range("K58").Value = 1244.875
$ws.Range("L58").Value = 5615
$ws.Range("M58").Value = -1041.875
$ws.Range("N58").Value = -6021
$ws.Range("H101").Value = 47328.5
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 47328.5
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 47328.5
$ws.Range("N101").Value = -53818.5
$ws.Range("H105").Value = 2144.0833
$ws.Range("I105").Value = 1055
$ws.Range("J105").Value = 2688.625
$ws.Range("K105").Value = 1055
$ws.Range("L105").Value = 2688.625
$ws.Range("M105").Value = 692
$ws.Range("N105").Value = -6182.625
$ws.Range("H107").Value = 453.0909
$ws.Range("I107").Value = 438.75
$ws.Range("J107").Value = 461.2857
$ws.Range("K107").Value = 438.75
$ws.Range("L107").Value = 461.2857
$ws.Range("M107").Value = 1481.25
$ws.Range("N107").Value = -4301.2857
$ws.Range("H136").Value = 3284.2666
$ws.Range("I136").Value = 1244.875
$ws.Range("J136").Value = 5615
$ws.Range("K136").Value = 3734.625
$ws.Range("L136").Value = 16845
$ws.Range("M136").Value = -1184.625
$ws.Range("N136").Value = -21945

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 28.90909
$ws.Range("I12").Value = 15.2
$ws.Range("J12").Value = 32.941177
$ws.Range("K12").Value = 45.59999999999999
$ws.Range("L12").Value = 98.823531
$ws.Range("M12").Value = 127.4
$ws.Range("N12").Value = -444.823531
$ws.Range("H23").Value = 89.25
$ws.Range("I23").Value = 60.142857
$ws.Range("J23").Value = 130
$ws.Range("K23").Value = 180.428571
$ws.Range("L23").Value = 390
$ws.Range("M23").Value = 54.57142899999999
$ws.Range("N23").Value = -860
$ws.Range("H33").Value = 108.333336
$ws.Range("I33").Value = 108.333336
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 650.000016
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -367.000016
$ws.Range("N33").ClearContents()
$ws.Range("H38").Value = 178.38461
$ws.Range("I38").Value = 416
$ws.Range("J38").Value = 72.77778000000001
$ws.Range("K38").Value = 1248
$ws.Range("L38").Value = 218.33334
$ws.Range("M38").Value = -901
$ws.Range("N38").Value = -912.33334
$ws.Range("H80").Value = 4722.222
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 5083.3335
$ws.Range("K80").Value = 12000
$ws.Range("L80").Value = 15250.0005
$ws.Range("M80").Value = -11064
$ws.Range("N80").Value = -17122.0005
$ws.Range("H83").Value = 4722.222
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 5083.3335
$ws.Range("K83").Value = 36000
$ws.Range("L83").Value = 45750.0015
$ws.Range("M83").Value = -31320
$ws.Range("N83").Value = -55110.0015
$ws.Range("H125").Value = 5420
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 5420
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 16260
$ws.Range("N125").Value = -26100

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1287.1305
$ws.Range("I97").Value = 1024
$ws.Range("J97").Value = 1780.5
$ws.Range("K97").Value = 1024
$ws.Range("L97").Value = 1780.5
$ws.Range("M97").Value = -528
$ws.Range("N97").Value = -2772.5
$ws.Range("H122").Value = 3705691
$ws.Range("I122").Value = 6062326.5
$ws.Range("J122").Value = 2406.1428
$ws.Range("K122").Value = 18186979.5
$ws.Range("L122").Value = 7218.428400000001
$ws.Range("M122").Value = -18184529.5
$ws.Range("N122").Value = -12118.4284
$ws.Range("H132").Value = 4601.1562
$ws.Range("I132").Value = 3535.625
$ws.Range("J132").Value = 5666.6875
$ws.Range("K132").Value = 10606.875
$ws.Range("L132").Value = 17000.0625
$ws.Range("M132").Value = -8076.875
$ws.Range("N132").Value = -22060.0625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1059.2
$ws.Range("I46").Value = 825
$ws.Range("J46").Value = 1144.3636
$ws.Range("K46").Value = 825
$ws.Range("L46").Value = 1144.3636
$ws.Range("M46").Value = -637
$ws.Range("N46").Value = -1520.3636
$ws.Range("H132").Value = 11913707
$ws.Range("I132").Value = 5873.9565
$ws.Range("J132").Value = 26328454
$ws.Range("K132").Value = 17621.8695
$ws.Range("L132").Value = 78985362
$ws.Range("M132").Value = -15091.8695
$ws.Range("N132").Value = -78990422
$ws.Range("H134").Value = 28308.666
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 28308.666
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 28308.666
$ws.Range("N134").Value = -38448.666

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2113.5881
$ws.Range("I132").Value = 1359.2858
$ws.Range("J132").Value = 2641.6
$ws.Range("K132").Value = 4077.8574
$ws.Range("L132").Value = 7924.799999999999
$ws.Range("M132").Value = -1547.8574
$ws.Range("N132").Value = -12984.8
$ws.Range("H133").Value = 39273
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 39273
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 39273
$ws.Range("N133").Value = -49393
